$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the added "Matières enseignés" column (E1)
$ws.Range("E1").Value = "Matières enseignés"

# Column widths for C, D and E (values chosen so the engine's internal
# pixel-rounded ColumnWidth->stored-width conversion lands as close as
# possible to the target stored widths 27.5703125 / 15.7109375 / 31.7109375)
$ws.Columns.Item(3).ColumnWidth = 26.6666666666667
$ws.Columns.Item(4).ColumnWidth = 14.8333333333333
$ws.Columns.Item(5).ColumnWidth = 30.8333333333333

# Move/restore the active selection to E6, matching the saved view state
$ws.Range("E6").Select()
